$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Range("A1").Value = 'Datos actualizados a 18 de Mayo de 2020 a las 19:35'

# Row 4
$ws.Range("B4").Value = 1534977
$ws.Range("C4").Value = 7313
$ws.Range("D4").Value = 347702
$ws.Range("E4").Value = 1095969
$ws.Range("G4").Value = 328
$ws.Range("H4").Value = 91306

# Row 8
$ws.Range("B8").Value = 245595
$ws.Range("C8").Value = 4515
$ws.Range("E8").Value = 135103
$ws.Range("G8").Value = 252
$ws.Range("H8").Value = 16370

# Row 12
$ws.Range("B12").Value = 150593
$ws.Range("C12").Value = 1158
$ws.Range("D12").Value = 111577
$ws.Range("E12").Value = 34845
$ws.Range("G12").Value = 31
$ws.Range("H12").Value = 4171

# Row 14
$ws.Range("B14").Value = 100340
$ws.Range("C14").Value = 4642
$ws.Range("D14").Value = 39231
$ws.Range("E14").Value = 57954

# Row 31
$ws.Range("B31").Value = 24200
$ws.Range("C31").Value = 88
$ws.Range("E31").Value = 3183
$ws.Range("G31").Value = 4
$ws.Range("H31").Value = 1547

# Row 37
$ws.Range("E37").Value = 5986
$ws.Range("G37").Value = 13
$ws.Range("H37").Value = 1120

# Row 44
$ws.Range("A44").Value = 'Egipto'
$ws.Range("B44").Value = 12764
$ws.Range("C44").Value = 535
$ws.Range("D44").Value = 3440
$ws.Range("E44").Value = 8679
$ws.Range("G44").Value = 15
$ws.Range("H44").Value = 645

# Row 45
$ws.Range("A45").Value = 'Republica Dominicana'
$ws.Range("B45").Value = 12725
$ws.Range("C45").Value = 411
$ws.Range("D45").Value = 6613
$ws.Range("E45").Value = 5678
$ws.Range("G45").Value = 6
$ws.Range("H45").Value = 434

# Row 46
$ws.Range("A46").Value = 'Filipinas'
$ws.Range("B46").Value = 12718
$ws.Range("C46").Value = 205
$ws.Range("D46").Value = 2729
$ws.Range("E46").Value = 9158
$ws.Range("G46").Value = 7
$ws.Range("H46").Value = 831

# Row 68
$ws.Range("B68").Value = 3947
$ws.Range("C68").Value = 2
$ws.Range("D68").Value = 3715
$ws.Range("E68").Value = 125

# Row 105
$ws.Range("A105").Value = 'Sri Lanka'
$ws.Range("B105").Value = 991
$ws.Range("C105").Value = 10
$ws.Range("D105").Value = 559
$ws.Range("E105").Value = 423
$ws.Range("H105").Value = 9

# Row 106
$ws.Range("A106").Value = 'Guinea-Bisau'
$ws.Range("B106").Value = 990
$ws.Range("C106").Value = 0
$ws.Range("D106").Value = 26
$ws.Range("E106").Value = 960
$ws.Range("H106").Value = 4

# Row 115
$ws.Range("B115").Value = 788
$ws.Range("C115").Value = 2
$ws.Range("D115").Value = 219
$ws.Range("E115").Value = 558
